$d = $word.ActiveDocument

$replacements = @(
    @("52×19=988", "35×28=980"),
    @("33×91=3003", "91×41=3731"),
    @("86×62=5332", "48×78=3744"),
    @("57×59=3363", "17×34=578"),
    @("50×85=4250", "43×15=645"),
    @("17×30=510", "96×82=7872"),
    @("12×88=1056", "84×84=7056"),
    @("93×13=1209", "97×82=7954"),
    @("15×39=585", "88×76=6688"),
    @("20×26=520", "84×28=2352"),
    @("31×98=3038", "37×39=1443"),
    @("54×56=3024", "88×95=8360"),
    @("93×35=3255", "22×78=1716"),
    @("80×11=880", "57×62=3534"),
    @("84×34=2856", "30×92=2760"),
    @("27×23=621", "16×22=352"),
    @("47×98=4606", "67×74=4958"),
    @("15×95=1425", "82×33=2706"),
    @("62×53=3286", "55×38=2090"),
    @("63×59=3717", "67×38=2546"),
    @("37×59=2183", "13×97=1261"),
    @("26×77=2002", "22×15=330"),
    @("41×40=1640", "18×89=1602"),
    @("25×73=1825", "97×48=4656"),
    @("55×77=4235", "20×69=1380")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
